$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 122, shifting existing rows 122:225 down to 123:226
$ws.Rows(122).Insert()

# Populate the newly inserted row 122 with the new data entry
$ws.Range("A122").Value = 11
$ws.Range("B122").Value = "Vega Monumental Concepción"
$ws.Range("C122").Value = "Bíobío"
$ws.Range("D122").Value = 44818
$ws.Range("E122").Value = 8
$ws.Range("F122").Value = 100112040
$ws.Range("G122").Value = "Cilantro"
$ws.Range("H122").Value = "Sin especificar"
$ws.Range("I122").Value = "Primera"
$ws.Range("J122").Value = 90
$ws.Range("K122").Value = 5000
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = 5556
$ws.Range("N122").Value = "`$/caja 36 atados"
$ws.Range("O122").Value = "Región Metropolitana"
$ws.Range("P122").Value = 154
$ws.Range("Q122").Value = 36
$ws.Range("R122").Value = "Hortaliza"
